# Updates cryptos list: Price (D) and Volume(1h) (E) columns
# D-column values are written with a leading apostrophe + style reset
# to guarantee they persist as literal text (matching the original
# inlineStr cell type) even when the digits would otherwise parse as a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.717.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.39%  '
$ws.Range("D3").Value = "'3.119.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.20%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = "'532.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.09%  '
$ws.Range("D6").Value = "'138.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.19%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = "'3.118.43"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.21%  '
$ws.Range("E9").Value = '  +4.85%  '
$ws.Range("D10").Value = "'7.32"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.52%  '
$ws.Range("D11").Value = "'0.107"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.26%  '
$ws.Range("D12").Value = "'0.414"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.41%  '
$ws.Range("D13").Value = "'3.656.15"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.16%  '
$ws.Range("E14").Value = '  +1.36%  '
$ws.Range("D15").Value = "'25.44"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.46%  '
$ws.Range("E16").Value = '  +0.33%  '
$ws.Range("D17").Value = "'57.873.25"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.42%  '
$ws.Range("D18").Value = "'3.122.69"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.17%  '
$ws.Range("D19").Value = "'6.01"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.04%  '
$ws.Range("D20").Value = "'12.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.76%  '
$ws.Range("D21").Value = "'8.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'361.15"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.49%  '
$ws.Range("E23").Value = '  -0.18%  '
$ws.Range("D24").Value = "'69.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.02%  '
$ws.Range("D25").Value = "'0.504"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("E26").Value = '  -0.56%  '
$ws.Range("E27").Value = '  +0.31%  '
$ws.Range("D28").Value = "'0.0₃0874"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.21%  '
$ws.Range("D29").Value = "'7.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.72%  '
$ws.Range("E30").Value = '  -0.29%  '
$ws.Range("E31").Value = '  -0.21%  '
$ws.Range("D32").Value = "'21.36"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.21%  '
$ws.Range("D33").Value = "'5.13"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.38%  '
$ws.Range("E34").Value = '  -2.52%  '
$ws.Range("D35").Value = "'158.31"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.08%  '
$ws.Range("D36").Value = "'6.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.62%  '
$ws.Range("D37").Value = "'25.82"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.29%  '
$ws.Range("E38").Value = '  +2.00%  '
$ws.Range("D39").Value = "'1.65"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.54%  '
$ws.Range("D40").Value = "'0.0672"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.06%  '
$ws.Range("D41").Value = "'2.484.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.95%  '
$ws.Range("D42").Value = "'0.697"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.44%  '
$ws.Range("D43").Value = "'3.99"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.07%  '
$ws.Range("D44").Value = "'37.65"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.14%  '
$ws.Range("D45").Value = "'3.162.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.29%  '
$ws.Range("E46").Value = '  -0.04%  '
$ws.Range("D47").Value = "'0.0268"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.21%  '
$ws.Range("D48").Value = "'0.987"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.89%  '
$ws.Range("D49").Value = "'6.06"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.44%  '
$ws.Range("D50").Value = "'19.77"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.70%  '
$ws.Range("D51").Value = "'0.741"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.88%  '
